$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "66.951.15"
$ws.Range("E2").Value = "  -3.67%  "

$ws.Range("D3").Value = "3.341.94"
$ws.Range("E3").Value = "  -0.72%  "

$ws.Range("E4").Value = "  -0.01%  "

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "574.51"
$ws.Range("E5").Value = "  -3.20%  "

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "182.55"
$ws.Range("E6").Value = "  -4.94%  "

$ws.Range("E7").Value = "  -0.01%  "

$ws.Range("E8").Value = "  -1.64%  "

$ws.Range("E9").Value = "  -3.15%  "

$ws.Range("E10").Value = "  -1.51%  "

$ws.Range("E11").Value = "  -4.43%  "

$ws.Range("D12").Value = "3.923.72"
$ws.Range("E12").Value = "  -0.76%  "

$ws.Range("E13").Value = "  -1.22%  "

$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "27.25"
$ws.Range("E14").Value = "  -4.88%  "

$ws.Range("D15").Value = "66.974.65"
$ws.Range("E15").Value = "  -3.65%  "

$ws.Range("E16").Value = "  -2.21%  "

$ws.Range("D17").Value = "3.343.25"
$ws.Range("E17").Value = "  -0.67%  "

$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "436.89"
$ws.Range("E18").Value = "  -2.95%  "

$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "13.71"
$ws.Range("E19").Value = "  -1.08%  "

$ws.Range("E20").Value = "  -2.52%  "

$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "7.66"

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "73.87"
$ws.Range("E22").Value = "  -0.84%  "

$ws.Range("E23").Value = "  +0.13%  "

$ws.Range("E25").Value = "  -2.48%  "

$ws.Range("E26").Value = "  -0.92%  "

$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "9.14"
$ws.Range("E27").Value = "  -4.63%  "

$ws.Range("E28").Value = "  +0.12%  "

$ws.Range("E29").Value = "  -1.53%  "

$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "22.88"
$ws.Range("E30").Value = "  -1.80%  "

$ws.Range("E31").Value = "  -4.70%  "

$ws.Range("B32").Value = "USDe"
$ws.Range("C32").Value = "https://coinranking.com/coin/exbfr2U-0+usde-usde"
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "0.999"
$ws.Range("E32").Value = "  +0.05%  "

$ws.Range("B33").Value = "Fetch.AI"
$ws.Range("C33").Value = "https://coinranking.com/coin/AWma-WzFHmKVQ+fetchai-fet"
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "1.24"
$ws.Range("E33").Value = "  -4.05%  "

$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "6.83"
$ws.Range("E34").Value = "  -2.91%  "

$ws.Range("E35").Value = "  -1.54%  "

$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "161.13"
$ws.Range("E36").Value = "  -2.43%  "

$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "27.72"
$ws.Range("E37").Value = "  +1.78%  "

$ws.Range("E38").Value = "  -4.80%  "

$ws.Range("D39").Value = "2.839.35"
$ws.Range("E39").Value = "  +3.60%  "

$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "0.793"
$ws.Range("E40").Value = "  -2.98%  "

$ws.Range("E41").Value = "  -3.49%  "

$ws.Range("E42").Value = "  -4.89%  "

$ws.Range("E43").Value = "  -1.84%  "

$ws.Range("E44").Value = "  -0.72%  "

$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "24.71"
$ws.Range("E45").Value = "  -4.42%  "

$ws.Range("E46").Value = "  -6.27%  "

$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "326.32"
$ws.Range("E47").Value = "  -5.03%  "

$ws.Range("E48").Value = "  -4.22%  "

$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "31.56"
$ws.Range("E49").Value = "  -4.72%  "

$ws.Range("E50").Value = "  -4.91%  "

$ws.Range("E51").Value = "  -2.56%  "
